# Fixar o banco de dados de autoria e mencoes no ano de análise
# Update the "categoria_mencao" values in sheets "max-arrecad" and
# "tx-sucesso" to reflect the re-computed ranking (values in column B are
# unchanged; only the category labels in column A are corrected).

$wb = $excel.ActiveWorkbook

# --- Sheet "max-arrecad": rows 2-10 ---
$ws4 = $wb.Worksheets.Item("max-arrecad")
$ws4.Range("A2").Value = "midia_independente"
$ws4.Range("A3").Value = "disputa"
$ws4.Range("A4").Value = "herois"
$ws4.Range("A5").Value = "terror"
$ws4.Range("A6").Value = "politica"
$ws4.Range("A7").Value = "religiosidade"
$ws4.Range("A8").Value = "jogos"
$ws4.Range("A9").Value = "erotismo"
$ws4.Range("A10").Value = "humor"

# --- Sheet "max-arrecad": rows 18-19 ---
$ws4.Range("A18").Value = "hqmix"
$ws4.Range("A19").Value = "questoes_genero"

# --- Sheet "tx-sucesso": rows 5-6 ---
$ws5 = $wb.Worksheets.Item("tx-sucesso")
$ws5.Range("A5").Value = "saloes_humor"
$ws5.Range("A6").Value = "questoes_genero"
